$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "CMS"

$headers = @(
    "Contact_ID`n",
    "Contact_Date`n",
    "Contact_Type_Code",
    "Contact_Type_Desc",
    "Contact_Staff_Name",
    "Contact_Staff_Key",
    "Contact_Staff_Grade",
    "Contact_Team_Key",
    "Contact_Provider_Code",
    "OM_Name`n",
    "OM_Key`n",
    "OM_Grade`n",
    "OM_Team_Key`n",
    "OM_Provider_Code`n"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
}

$headerRange = $newSheet.Range("A1:N1")
$font = $headerRange.Font
$font.Name = "Arial"
$font.Size = 9
$font.Color = 3355443
$headerRange.Interior.Color = 16777215
$headerRange.HorizontalAlignment = -4131
$newSheet.Rows.Item(1).RowHeight = 23.25

$newSheet.Range("A1").Select()
